$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs","Mmp2","Sdc2","ECs",3,1,10.29746733333333,30.892402,0.03236282443656328,0.03236282443656328,3,1,1.826566,5.479698,0.02795372904983374,0.02795372904983374,18.80900371717733,169.281033454596,0.0009046616255870281,0.0009046616255870283),
    @("ECs","Mmp2","Sdc2","FAPs",3,1,10.29746733333333,30.892402,0.03236282443656328,0.03236282443656328,3,1,44.29005966666667,132.870179,0.6778141756295529,0.6778141756295529,456.0754426088843,4104.678983479958,0.02193598116651309,0.02193598116651309),
    @("ECs","Mmp2","Sdc2","MuSCs",3,1,10.29746733333333,30.892402,0.03236282443656328,0.03236282443656328,3,1,19.10886933333333,57.326608,0.2924417490485847,0.2924417490485847,196.7729577369351,1770.956619632416,0.009464240982380844,0.009464240982380844),
    @("ECs","Mmp2","Sdc2","Resolving-Mac",3,1,10.29746733333333,30.892402,0.03236282443656328,0.03236282443656328,2,0.6666666666666666,0.1169856666666667,0.350957,0.001790346272028586,0.001790346272028586,1.204656080968222,10.841904728714,0.0000579406620823167,0.0000579406620823167),
    @("FAPs","Mmp2","Sdc2","ECs",3,1,301.5732116666666,904.7196349999999,0.9477826525699299,0.9477826525699298,3,1,1.826566,5.479698,0.02795372904983374,0.02795372904983374,550.8433749411365,4957.590374470229,0.02649405946807253,0.02649405946807253),
    @("FAPs","Mmp2","Sdc2","FAPs",3,1,301.5732116666666,904.7196349999999,0.9477826525699299,0.9477826525699298,3,1,44.29005966666667,132.870179,0.6778141756295529,0.6778141756295529,13356.69553858496,120210.2598472647,0.642420517327678,0.6424205173276779),
    @("FAPs","Mmp2","Sdc2","MuSCs",3,1,301.5732116666666,904.7196349999999,0.9477826525699299,0.9477826525699298,3,1,19.10886933333333,57.326608,0.2924417490485847,0.2924417490485847,5762.723096172009,51864.50786554808,0.2771712166354574,0.2771712166354574),
    @("FAPs","Mmp2","Sdc2","Resolving-Mac",3,1,301.5732116666666,904.7196349999999,0.9477826525699299,0.9477826525699298,2,0.6666666666666666,0.1169856666666667,0.350957,0.001790346272028586,0.001790346272028586,35.27974321563277,317.5176889406949,0.001696859138721939,0.001696859138721939),
    @("Inflammatory-Mac","Mmp2","Sdc2","ECs",3,1,1.128766333333333,3.386299,0.003547480705019629,0.003547480705019629,3,1,1.826566,5.479698,0.02795372904983374,0.02795372904983374,2.061766206411333,18.555895857702,0.00009916531443763187,0.00009916531443763189),
    @("Inflammatory-Mac","Mmp2","Sdc2","FAPs",3,1,1.128766333333333,3.386299,0.003547480705019629,0.003547480705019629,3,1,44.29005966666667,132.870179,0.6778141756295529,0.6778141756295529,49.9931282530579,449.9381542775211,0.002404532709634625,0.002404532709634625),
    @("Inflammatory-Mac","Mmp2","Sdc2","MuSCs",3,1,1.128766333333333,3.386299,0.003547480705019629,0.003547480705019629,3,1,19.10886933333333,57.326608,0.2924417490485847,0.2924417490485847,21.56944837153245,194.125035343792,0.001037431462092047,0.001037431462092047),
    @("Inflammatory-Mac","Mmp2","Sdc2","Resolving-Mac",3,1,1.128766333333333,3.386299,0.003547480705019629,0.003547480705019629,2,0.6666666666666666,0.1169856666666667,0.350957,0.001790346272028586,0.001790346272028586,0.1320494820158889,1.188445338143,0.000006351218855325233,0.000006351218855325233),
    @("MuSCs","Mmp2","Sdc2","ECs",3,1,3.752174,11.256522,0.01179231207894783,0.01179231207894783,3,1,1.826566,5.479698,0.02795372904983374,0.02795372904983374,6.853593454484,61.682341090356,0.0003296390967259893,0.0003296390967259893),
    @("MuSCs","Mmp2","Sdc2","FAPs",3,1,3.752174,11.256522,0.01179231207894783,0.01179231207894783,3,1,44.29005966666667,132.870179,0.6778141756295529,0.6778141756295529,166.1840103397153,1495.656093057438,0.007992996290558444,0.007992996290558442),
    @("MuSCs","Mmp2","Sdc2","MuSCs",3,1,3.752174,11.256522,0.01179231207894783,0.01179231207894783,3,1,19.10886933333333,57.326608,0.2924417490485847,0.2924417490485847,71.69980268193068,645.298224137376,0.003448564369694256,0.003448564369694256),
    @("MuSCs","Mmp2","Sdc2","Resolving-Mac",3,1,3.752174,11.256522,0.01179231207894783,0.01179231207894783,2,0.6666666666666666,0.1169856666666667,0.350957,0.001790346272028586,0.001790346272028586,0.4389505768393333,3.950555191554,0.00002111232196914192,0.00002111232196914192),
    @("Resolving-Mac","Mmp2","Sdc2","ECs",3,1,1.436533666666667,4.309601000000001,0.004514730209539471,0.00451473020953947,3,1,1.826566,5.479698,0.02795372904983374,0.02795372904983374,2.623923553388667,23.615311980498,0.0001262035450105655,0.0001262035450105655),
    @("Resolving-Mac","Mmp2","Sdc2","FAPs",3,1,1.436533666666667,4.309601000000001,0.004514730209539471,0.00451473020953947,3,1,44.29005966666667,132.870179,0.6778141756295529,0.6778141756295529,63.62416180984212,572.6174562885791,0.003060148135168835,0.003060148135168834),
    @("Resolving-Mac","Mmp2","Sdc2","MuSCs",3,1,1.436533666666667,4.309601000000001,0.004514730209539471,0.00451473020953947,3,1,19.10886933333333,57.326608,0.2924417490485847,0.2924417490485847,27.45053412926756,247.054807163408,0.001320295598960206,0.001320295598960206),
    @("Resolving-Mac","Mmp2","Sdc2","Resolving-Mac",3,1,1.436533666666667,4.309601000000001,0.004514730209539471,0.00451473020953947,2,0.6666666666666666,0.1169856666666667,0.350957,0.001790346272028586,0.001790346272028586,0.1680538486841111,1.512484638157,0.0000808293039986383,0.000008082930399863828)
)

$r = 2
foreach ($row in $data) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c = $c + 1
    }
    $r = $r + 1
}
